# Add a new "mailing list" block of 4 columns (Alain / Henri / Tony /
# Dulcinée) to "sheet1", inserted right before the existing JU:JV
# (email / trailing empty) columns. The insert pushes the existing
# JU:JV columns out to JY:JZ, and the 4 freshly inserted columns are
# populated by copying the pattern from the block immediately to their
# left (JQ:JT), matching the repeating 4-column cycle used across the
# rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# Insert 4 blank columns at JU:JX -- shifts old JU:JV (email + blank)
# to JY:JZ and widens the sheet's dimension accordingly.
$ws.Columns("JU:JX").Insert()

# Fill the newly inserted columns with the same repeating header/values
# pattern as the block to their left (now at JQ:JT).
$ws.Range("JQ1:JT9").Copy() | Out-Null
$ws.Range("JU1").PasteSpecial() | Out-Null
